# Updates the "cryptos" price/volume table (Sheet1) with freshly scraped
# values, as produced by the scheduled GitHub Actions job.
#
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# Price/volume cells are stored as plain text (e.g. "88.048.53",
# "  +7.72%  ") rather than numbers, so for any new price that looks like
# a plain number we force the cell to text format first - otherwise Excel
# would silently reinterpret the string as a numeric value and normalize/
# round it (e.g. "35.72" -> 35.7199999999999, "0.410" -> 0.41), which would
# not match the intended literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "87.935.98"
$ws.Range("E2").Value = "  +7.14%  "

Set-TextValue $ws.Range("D3") "3.365.95"
$ws.Range("E3").Value = "  +5.68%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws.Range("D5") "216.45"
$ws.Range("E5").Value = "  +2.20%  "

Set-TextValue $ws.Range("D6") "634.72"
$ws.Range("E6").Value = "  +1.74%  "

Set-TextValue $ws.Range("D7") "0.410"
$ws.Range("E7").Value = "  +41.86%  "

Set-TextValue $ws.Range("D8") "0.651"
$ws.Range("E8").Value = "  +11.27%  "

$ws.Range("E9").Value = "  -0.04%  "

Set-TextValue $ws.Range("D10") "3.362.27"
$ws.Range("E10").Value = "  +5.64%  "

Set-TextValue $ws.Range("D11") "0.602"
$ws.Range("E11").Value = "  +1.73%  "

Set-TextValue $ws.Range("D12") "0.0000274"
$ws.Range("E12").Value = "  +5.86%  "

Set-TextValue $ws.Range("D13") "35.72"
$ws.Range("E13").Value = "  +12.91%  "

Set-TextValue $ws.Range("D14") "0.168"
$ws.Range("E14").Value = "  +1.87%  "

Set-TextValue $ws.Range("D15") "3.971.11"
$ws.Range("E15").Value = "  +5.45%  "

Set-TextValue $ws.Range("D16") "5.38"
$ws.Range("E16").Value = "  +1.20%  "

Set-TextValue $ws.Range("D17") "87.579.24"
$ws.Range("E17").Value = "  +7.29%  "

Set-TextValue $ws.Range("D18") "3.360.33"
$ws.Range("E18").Value = "  +6.17%  "

Set-TextValue $ws.Range("D19") "14.55"
$ws.Range("E19").Value = "  +3.97%  "

Set-TextValue $ws.Range("D20") "9.40"
$ws.Range("E20").Value = "  +5.23%  "

Set-TextValue $ws.Range("D21") "447.09"
$ws.Range("E21").Value = "  +2.67%  "

Set-TextValue $ws.Range("D22") "2.98"
$ws.Range("E22").Value = "  -7.14%  "

Set-TextValue $ws.Range("D23") "5.48"
$ws.Range("E23").Value = "  +7.55%  "

Set-TextValue $ws.Range("D24") "7.25"
$ws.Range("E24").Value = "  -0.17%  "

Set-TextValue $ws.Range("D25") "5.39"
$ws.Range("E25").Value = "  +2.36%  "

Set-TextValue $ws.Range("D26") "12.40"
$ws.Range("E26").Value = "  +14.07%  "

Set-TextValue $ws.Range("D27") "3.547.85"
$ws.Range("E27").Value = "  +7.31%  "

Set-TextValue $ws.Range("D28") "79.16"
$ws.Range("E28").Value = "  +3.12%  "

Set-TextValue $ws.Range("D29") "0.0000137"
$ws.Range("E29").Value = "  +11.88%  "

Set-TextValue $ws.Range("D30") "0.997"
$ws.Range("E30").Value = "  -0.75%  "

Set-TextValue $ws.Range("D31") "0.177"
$ws.Range("E31").Value = "  +28.24%  "

Set-TextValue $ws.Range("D32") "9.13"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("E33").Value = "  +0.28%  "

Set-TextValue $ws.Range("D34") "564.61"
$ws.Range("E34").Value = "  -4.08%  "

Set-TextValue $ws.Range("D35") "1.50"
$ws.Range("E35").Value = "  -0.08%  "

Set-TextValue $ws.Range("D36") "2.05"
$ws.Range("E36").Value = "  +2.83%  "

Set-TextValue $ws.Range("D37") "7.14"
$ws.Range("E37").Value = "  +16.14%  "

Set-TextValue $ws.Range("D38") "0.139"
$ws.Range("E38").Value = "  -9.95%  "

Set-TextValue $ws.Range("D39") "23.35"
$ws.Range("E39").Value = "  +2.61%  "

Set-TextValue $ws.Range("D40") "0.422"
$ws.Range("E40").Value = "  +3.47%  "

Set-TextValue $ws.Range("D41") "21.87"
$ws.Range("E41").Value = "  +5.24%  "

$ws.Range("E42").Value = "  -0.02%  "

Set-TextValue $ws.Range("D43") "2.07"
$ws.Range("E43").Value = "  +0.04%  "

Set-TextValue $ws.Range("D44") "3.01"
$ws.Range("E44").Value = "  -1.84%  "

$ws.Range("E45").Value = "  -0.13%  "

Set-TextValue $ws.Range("D46") "157.01"
$ws.Range("E46").Value = "  -1.82%  "

Set-TextValue $ws.Range("D47") "183.40"
$ws.Range("E47").Value = "  -2.41%  "

Set-TextValue $ws.Range("D48") "1.39"
$ws.Range("E48").Value = "  +3.39%  "

Set-TextValue $ws.Range("D49") "45.92"
$ws.Range("E49").Value = "  +2.11%  "

Set-TextValue $ws.Range("D50") "4.39"
$ws.Range("E50").Value = "  +4.20%  "

# Row 51: ARBITRUM -> Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D51") "0.768"
$ws.Range("E51").Value = "  -0.94%  "
